# Fill every currently-blank cell in C2:K55 with a literal 0, leaving
# already-populated cells untouched, then move the active-cell selection
# from J5 to J9 (as captured in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 55

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 3; $c -le 11; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value() -eq $null) {
            $cell.Value = 0
        }
    }
}

$ws.Range("J9").Select()
